# Model Plans.xlsx edit:
# Add a new "Model Paddy 5" row to the crop model plan table, right after
# the existing "Model Paddy 4" row (and before "Model Wheat 1"), since the
# Paddy growing season now spans one additional month (sowing moved from
# 1-Jun to 1-Jul, harvest moved from 1-Sep to 1-Oct).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Model Paddy 4" currently lives in row 29 and "Model Wheat 1" in row 30.
# Insert a new blank row at 30, shifting "Model Wheat 1" (and everything
# below it) down by one.
$insertRowNum = 30
$ws.Rows.Item($insertRowNum).Insert(-4121)

# Copy the formatting of the row directly above (the last "Paddy" row) onto
# the freshly inserted row so it matches the look of the rest of the table.
$ws.Range("A29:D29").Copy()
$ws.Range("A30:D30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the row height used by the surrounding Paddy / Wheat rows.
$ws.Rows.Item($insertRowNum).RowHeight = $ws.Rows.Item(29).RowHeight

# Fill in the new row's values.
$ws.Cells.Item($insertRowNum, 1).Value2 = 29
$ws.Cells.Item($insertRowNum, 2).Value2 = "Model Paddy 5"
$ws.Cells.Item($insertRowNum, 3).Value2 = "Paddy Prediction - Month 5"
$ws.Cells.Item($insertRowNum, 4).Value2 = "1. Location, 2. Sowing Time Paddy(Farmer Provided), 3. Soil Nutrient (Farmer Provided), 4. Water Cycle, 5. Weather - 15 days prediction, 6. Weather Month 1, 7. Weather Month 2, 8. Wether Month 3, 9.Weather Month 4, 10. Weather Month 5"

# Renumber the S.No. column for every row that got pushed down by the insert.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = $insertRowNum + 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 1
    }
}

# Grow the "Table2" ListObject / AutoFilter to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B1:D41"))

# Reset the view back to the top-left of the sheet.
$ws.Range("B1").Select()
